# Unit8_Conditional Statements.pptx - rename leftover "Week3" footer label to
# "Unit8" (the deck was cloned from an earlier "Week3" deck; the footer slide
# number placeholder on the two CS1010-schedule exercise slides still read
# "Week3 - <n>" and needs to read "Unit8 - <n>").
#
# NB: the underlying commit also touched purely PowerPoint-internal metadata
# (ppt/revisionInfo.xml, ppt/changesInfos/*, and the cached text of the
# auto-updating "datetimeFigureOut" date field on the Notes Master, which
# PowerPoint itself stamps with the current date and which is not an
# editable/settable part of the document object model) - those are not
# user content and are out of scope for a COM automation script.

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.Name -eq "Slide Number Placeholder 8" -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Length -ge 6 -and $tr.Characters(1, 6).Text -eq "Week3 ") {
                $tr.Characters(1, 6).Text = "Unit8 "
            }
        }
    }
}
